$wb = $excel.ActiveWorkbook

# Update the header text in both sheets: "FACT / REMIS" -> "FACT/REMIS"
$wsPendientes = $wb.Worksheets.Item("Pendientes")
$wsFacturados = $wb.Worksheets.Item("Facturados")

$wsPendientes.Range("C1").Value = "FACT/REMIS"
$wsFacturados.Range("C1").Value = "FACT/REMIS"

# Update selections on each sheet
$wsPendientes.Range("C1").Select()
$wsFacturados.Range("C1").Select()

# Make Facturados the active sheet/tab
$wsFacturados.Activate()
